# 25th April 1st update
# Insert a new date column "25/04/2020" in its correct sorted position
# (right before the existing "26/03/2020" column, i.e. new column AZ),
# shifting the old AZ..BF columns to BA..BG. Also apply the handful of
# data corrections/additions that came with this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column at AZ (pushes old AZ:BF -> BA:BG)
$ws.Columns("AZ:AZ").Insert()

# New column header
$ws.Range("AZ1").Value = "25/04/2020"

# Make sure every data row gets an (empty) cell in the new column, matching
# the layout of the rest of the sheet, before filling in the real numbers.
$ws.Range("AZ2:AZ34").Style = "Normal"

# New data for 25/04/2020 (column AZ) on rows that have a reported count
$ws.Range("AZ4").Value = 61    # Andhra Pradesh
$ws.Range("AZ9").Value = 1     # Chhattisgarh
$ws.Range("AZ17").Value = 15   # Karnataka
$ws.Range("AZ28").Value = 25   # Rajasthan
$ws.Range("AZ34").Value = 57   # West Bengal

# Updated/corrected data for 24/04/2020 (column AX, unaffected by the insert)
$ws.Range("AX2").Value = $null   # blank state row - corrected to blank
$ws.Range("AX3").Value = 7       # Andaman and Nicobar Islands
$ws.Range("AX7").Value = 53      # Bihar
$ws.Range("AX10").Value = 138    # Delhi
$ws.Range("AX16").Value = 6      # Jharkhand
$ws.Range("AX28").Value = 70     # Rajasthan
$ws.Range("AX30").Value = 13     # Telangana
